$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.979605666666667
$ws.Range("N2").Value = 5.938817
$ws.Range("O2").Value = 0.05865520690928468
$ws.Range("P2").Value = 0.05865520690928468
$ws.Range("Q2").Value = 0.3517046611626667
$ws.Range("R2").Value = 3.165341950464
$ws.Range("S2").Value = 0.05865520690928468
$ws.Range("T2").Value = 0.05865520690928468

# Row 3
$ws.Range("M3").Value = 16.37791666666667
$ws.Range("N3").Value = 49.13375000000001
$ws.Range("O3").Value = 0.4852734597612733
$ws.Range("P3").Value = 0.4852734597612734
$ws.Range("Q3").Value = 2.909766186666667
$ws.Range("R3").Value = 26.18789568
$ws.Range("S3").Value = 0.4852734597612733
$ws.Range("T3").Value = 0.4852734597612734

# Row 4
$ws.Range("M4").Value = 15.39234866666667
$ws.Range("N4").Value = 46.177046
$ws.Range("O4").Value = 0.4560713333294419
$ws.Range("P4").Value = 0.4560713333294419
$ws.Range("Q4").Value = 2.734666233514667
$ws.Range("R4").Value = 24.611996101632
$ws.Range("S4").Value = 0.4560713333294419
$ws.Range("T4").Value = 0.4560713333294419
